$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.113595128059387
$ws.Range("B1").Value = 2.250239372253418
$ws.Range("C1").Value = 10.22958755493164
$ws.Range("D1").Value = 1.582264542579651
$ws.Range("E1").Value = 1.292576909065247
